$d = $word.ActiveDocument
$t = $d.Tables(1)

# ---------------------------------------------------------------
# 1. Insert two narrow "spacer" columns: before former column 2
#    (Course) and before former column 3 (Major, now at index 4
#    after the first insert). Result layout:
#    Name | spacer | Course | spacer | Major
# ---------------------------------------------------------------
$t.Columns.Add($t.Columns(2))
$t.Columns.Add($t.Columns(4))

# ---------------------------------------------------------------
# 2. Set the final column widths (dxa). Word.Column.Width is in
#    points, so divide the target dxa (twentieths of a point) by 20.
# ---------------------------------------------------------------
$t.Columns(1).Width = 3119 / 20.0
$t.Columns(2).Width = 283 / 20.0
$t.Columns(3).Width = 3402 / 20.0
$t.Columns(4).Width = 284 / 20.0
$t.Columns(5).Width = 3702 / 20.0

# ---------------------------------------------------------------
# 3. Center-align the (empty) paragraph inside every new spacer
#    cell, for every row.
# ---------------------------------------------------------------
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $t.Cell($r, 2).Range.ParagraphFormat.Alignment = 1
    $t.Cell($r, 4).Range.ParagraphFormat.Alignment = 1
}

# ---------------------------------------------------------------
# 4. Table-level borders: explicit "none" on every side.
# ---------------------------------------------------------------
$t.Borders.OutsideLineWidth = 0
$t.Borders.InsideLineWidth = 0
$t.Borders.Enable = 0

# ---------------------------------------------------------------
# 5. Per-cell borders on the data columns (1, 3, 5):
#      row 2 (first data row)      -> bottom only
#      rows 3-6 (remaining rows)   -> top + bottom
#    Using Borders.Item(<wdBorderXxx>) scopes the change to a
#    single edge instead of the whole cell.
# ---------------------------------------------------------------
$dataCols = 1, 3, 5

foreach ($c in $dataCols) {
    $cell = $t.Cell(2, $c)
    $bottom = $cell.Borders.Item(-3)
    $bottom.LineStyle = 1
    $bottom.LineWidth = 2
    $bottom.ColorIndex = 0
    $cell.Borders.DistanceFromBottom = 0
}

for ($r = 3; $r -le 6; $r++) {
    foreach ($c in $dataCols) {
        $cell = $t.Cell($r, $c)

        $top = $cell.Borders.Item(-1)
        $top.LineStyle = 1
        $top.LineWidth = 2
        $top.ColorIndex = 0
        $cell.Borders.DistanceFromTop = 0

        $bottom = $cell.Borders.Item(-3)
        $bottom.LineStyle = 1
        $bottom.LineWidth = 2
        $bottom.ColorIndex = 0
        $cell.Borders.DistanceFromBottom = 0
    }
}
